# Add 'ongkosKirimBeli' to the Produk (product) sheet and make it the
# active sheet/tab, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Produk")

# Insert a new blank column before column F (pushes modifiedDate .. jumlahTukar
# one column to the right: F->G, G->H, ... O->P).
$ws.Range("F1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 6).Value = "ongkosKirimBeli"

# Give the new column a sensible width (close to the authored best-fit width).
$ws.Columns.Item(6).ColumnWidth = 14.57

# Make "Produk" the active sheet/tab and select F2, mirroring the author's
# recorded view state (moves tabSelected from the old active sheet to this
# one, and updates the workbook's activeTab).
$ws.Activate()
$ws.Range("F2").Select()
